{"js": "// Append \". / Insert / Lookup / Delete / Search\" complexity list to the\n// end of the \"Hash tables\" section, move the `_GoBack` bookmark to the\n// new end of the document, and add a trailing empty paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The paragraph that currently ends the document:\n// \"...guardan en memoria un valor en base a una llave dada\"\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// 1) Close that sentence with a period (new run).\nlastParagraph.insertText(\".\", Word.InsertLocation.end);\nawait context.sync();\n\n// 2) The old `_GoBack` bookmark sits right at the (old) end of the\n// document; remove it here so it can be re-inserted at the new end\n// once the new paragraphs below are in place (Word always keeps\n// `_GoBack` pointing at the most-recently-edited spot).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst W_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nfunction wrapPackage(bodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    \"<w:document \" +\n    W_NS +\n    \">\" +\n    \"<w:body>\" +\n    bodyXml +\n    \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\n// Builds a paragraph like: <w:p><proofErr spellStart/><run word/><proofErr\n// spellEnd/><run \": \"/><proofErr gramStart/><run \"O(\"/><proofErr\n// gramEnd/><run \"1)\"/></w:p> \u2014 matching how Word's proofer marks these\n// runs up as you type them. The last one (\"Search\") is only spell-checked\n// (no gramStart/gramEnd split), matching the source edit exactly.\nfunction complexityParagraphXml(word, splitGrammar) {\n  const tail = splitGrammar\n    ? '<w:r><w:t xml:space=\"preserve\">: </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      \"<w:r><w:t>O(</w:t></w:r>\" +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      \"<w:r><w:t>1)</w:t></w:r>\"\n    : \"<w:r><w:t>: O(1)</w:t></w:r>\";\n  return (\n    \"<w:p>\" +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>\" +\n    word +\n    \"</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    tail +\n    \"</w:p>\"\n  );\n}\n\n// 3) Insert the four new paragraphs (Insert/Lookup/Delete/Search), each\n// \"O(1)\"), right after the paragraph we just closed with a period.\nconst newParagraphsXml =\n  complexityParagraphXml(\"Insert\", true) +\n  complexityParagraphXml(\"Lookup\", true) +\n  complexityParagraphXml(\"Delete\", true) +\n  complexityParagraphXml(\"Search\", false);\n\nconst insertionRange = lastParagraph.getRange(\"End\");\ninsertionRange.insertOoxml(wrapPackage(newParagraphsXml), Word.InsertLocation.after);\nawait context.sync();\n\n// 4) Re-create `_GoBack` collapsed at the new end of the body (i.e. right\n// after \"Search: O(1)\", before the trailing empty paragraph we add next).\nconst newEndRange = context.document.body.getRange(\"End\");\nnewEndRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 5) Add one trailing empty paragraph after \"Search: O(1)\".\nconst paragraphsAfter = context.document.body.paragraphs;\nparagraphsAfter.load(\"items\");\nawait context.sync();\nconst searchParagraph = paragraphsAfter.items[paragraphsAfter.items.length - 1];\nsearchParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Append the \"Insert/Lookup/Delete/Search: O(1)\" complexity list to the\n# end of the \"Hash tables\" section, move the `_GoBack` bookmark to the\n# new end of the document, and add a trailing empty paragraph.\n\n$d = $word.ActiveDocument\n\n# 0. Remove the existing `_GoBack` bookmark (if any); Word always keeps\n#    it collapsed at the most-recently-edited spot, so it will be\n#    re-created at the new end of the document once the new\n#    paragraphs below are in place.\ntry {\n    $oldBookmark = $d.Bookmarks.Item(\"_GoBack\")\n    $oldBookmark.Delete()\n} catch {\n}\n\n# 1. Close the current last paragraph's sentence with a period (own run).\n$endPos = $d.Content.End - 1\n$r = $d.Range($endPos, $endPos)\n$r.InsertAfter(\".\")\n\n# 2. Insert the four new paragraphs (each \"<Word>: O(1)\"), reproducing\n#    the spell/grammar-check proofing marks Word leaves behind as you\n#    type (\"Insert\"/\"Lookup\"/\"Delete\"/\"Search\" are flagged spellStart/\n#    spellEnd, and \"O(\" is flagged gramStart/gramEnd on every paragraph\n#    except the last one).\n$newParagraphsXml = (\n    '<w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Insert</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">: </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/><w:r><w:t>O(</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/><w:r><w:t>1)</w:t></w:r></w:p>' +\n    '<w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Lookup</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">: </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/><w:r><w:t>O(</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/><w:r><w:t>1)</w:t></w:r></w:p>' +\n    '<w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Delete</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">: </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/><w:r><w:t>O(</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/><w:r><w:t>1)</w:t></w:r></w:p>' +\n    '<w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Search</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/><w:r><w:t>: O(1)</w:t></w:r></w:p>'\n)\n\n$packageXml = (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + $newParagraphsXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n)\n\n$endPos = $d.Content.End - 1\n$r = $d.Range($endPos, $endPos)\n$null = $r.InsertXML($packageXml)\n\n# 3. Re-create `_GoBack`, collapsed right after \"Search: O(1)\".\n#    Bookmarks.Add() on a zero-width range touching a paragraph end is\n#    unreliable in this host, so: append a 1-character placeholder,\n#    wrap the bookmark around it (non-collapsed range), then clear the\n#    placeholder's text; the bookmark collapses in place, exactly where\n#    Word leaves `_GoBack` after the last edit.\n$endPos = $d.Content.End - 1\n$r = $d.Range($endPos, $endPos)\n$r.InsertAfter(\"X\")\n\n$endPos = $d.Content.End - 1\n$placeholderRange = $d.Range($endPos - 1, $endPos)\n$d.Bookmarks.Add(\"_GoBack\", $placeholderRange)\n\n$newBookmark = $d.Bookmarks.Item(\"_GoBack\")\n$newBookmark.Range.Text = \"\"\n\n# 4. Add one trailing empty paragraph after \"Search: O(1)\".\n$endPos = $d.Content.End - 1\n$r = $d.Range($endPos, $endPos)\n$r.InsertParagraphAfter()\n"}
